$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.505.01'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.902.41'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.75%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.69'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.68'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -6.73%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -1.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.906.79'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.81%  '
$ws.Range("E10").Value = '  -2.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.95'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.27%  '
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.411.17'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.80%  '
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.547.19'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.61'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.907.63'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.79%  '
$ws.Range("E18").Value = '  -3.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.05'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.70'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '363.67'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -6.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.62'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.14'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.024.99'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.453'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.84%  '
$ws.Range("E27").Value = '  -3.62%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.81'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -6.93%  '
$ws.Range("E30").Value = '  -7.65%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.68'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.65'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.10'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -7.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.37'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.58'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -6.66%  '
$ws.Range("E37").Value = '  -5.84%  '
$ws.Range("E38").Value = '  -5.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.91'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.91%  '
$ws.Range("E40").Value = '  -4.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.331.40'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.07%  '
$ws.Range("E42").Value = '  -5.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.645'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.71'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0576'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.997'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("E48").Value = '  -4.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0934'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.91%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.33'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.44'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.92%  '
